$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 1951.5385
$ws.Range("I55").Value = 1055.6
$ws.Range("J55").Value = 2511.5
$ws.Range("K55").Value = 1055.6
$ws.Range("L55").Value = 2511.5
$ws.Range("M55").Value = -841.5999999999999
$ws.Range("N55").Value = -2939.5

# Row 124
$ws.Range("H124").Value = 40000
$ws.Range("J124").Value = 40000
$ws.Range("L124").Value = 40000
$ws.Range("N124").Value = -49820

# Row 137
$ws.Range("H137").Value = 19897.83
$ws.Range("I137").Value = 37765.777
$ws.Range("J137").Value = 1342.6538
$ws.Range("K137").Value = 113297.331
$ws.Range("L137").Value = 4027.9614
$ws.Range("M137").Value = -110747.331
$ws.Range("N137").Value = -9127.9614

$ws = $wb.Worksheets.Item("ARM")
# Row 49
$ws.Range("H49").Value = 44770
$ws.Range("J49").Value = 44770
$ws.Range("L49").Value = 44770
$ws.Range("N49").Value = -45290

# Row 74
$ws.Range("H74").Value = 92473.41
$ws.Range("I74").Value = 119396.94
$ws.Range("J74").Value = 933.4
$ws.Range("K74").Value = 119396.94
$ws.Range("L74").Value = 933.4
$ws.Range("M74").Value = -118522.94
$ws.Range("N74").Value = -2681.4

# Row 77
$ws.Range("H77").Value = 92473.41
$ws.Range("I77").Value = 119396.94
$ws.Range("J77").Value = 933.4
$ws.Range("K77").Value = 596984.7
$ws.Range("L77").Value = 4667
$ws.Range("M77").Value = -592616.7
$ws.Range("N77").Value = -13403

# Row 86
$ws.Range("H86").Value = 166687500
$ws.Range("J86").Value = 166687500
$ws.Range("L86").Value = 166687500
$ws.Range("N86").Value = -166689872

# Row 89
$ws.Range("H89").Value = 166687500
$ws.Range("J89").Value = 166687500
$ws.Range("L89").Value = 500062500
$ws.Range("N89").Value = -500074356

# Row 132
$ws.Range("H132").Value = 5340474
$ws.Range("I132").Value = 8501895
$ws.Range("J132").Value = 1125245
$ws.Range("K132").Value = 25505685
$ws.Range("L132").Value = 3375735
$ws.Range("M132").Value = -25503155
$ws.Range("N132").Value = -3380795

$ws = $wb.Worksheets.Item("BSM")
# Row 14
$ws.Range("H14").Value = 42205.4
$ws.Range("J14").Value = 42205.4
$ws.Range("L14").Value = 42205.4
$ws.Range("N14").Value = -42549.4

# Row 134
$ws.Range("H134").Value = 54473.9
$ws.Range("I134").Value = 2580.3333
$ws.Range("K134").Value = 7740.999899999999
$ws.Range("M134").Value = -5205.999899999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 26316.238
$ws.Range("I31").Value = 44453.348
$ws.Range("J31").Value = 4360.7896
$ws.Range("K31").Value = 44453.348
$ws.Range("L31").Value = 4360.7896
$ws.Range("M31").Value = -44158.348
$ws.Range("N31").Value = -4950.7896

# Row 34
$ws.Range("H34").Value = 26316.238
$ws.Range("I34").Value = 44453.348
$ws.Range("J34").Value = 4360.7896
$ws.Range("K34").Value = 44453.348
$ws.Range("L34").Value = 4360.7896
$ws.Range("M34").Value = -44251.348
$ws.Range("N34").Value = -4764.7896

# Row 132
$ws.Range("H132").Value = 2822.9546
$ws.Range("I132").Value = 1379.1428
$ws.Range("J132").Value = 5349.625
$ws.Range("K132").Value = 4137.428400000001
$ws.Range("L132").Value = 16048.875
$ws.Range("M132").Value = -1607.428400000001
$ws.Range("N132").Value = -21108.875

# Row 134
$ws.Range("H134").Value = 10870927
$ws.Range("I134").Value = 1086.3715
$ws.Range("K134").Value = 3259.1145
$ws.Range("M134").Value = -724.1144999999997

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 2337.875
$ws.Range("I33").Value = 528.1429000000001
$ws.Range("J33").Value = 15006
$ws.Range("K33").Value = 3168.8574
$ws.Range("L33").Value = 90036
$ws.Range("M33").Value = -2885.8574
$ws.Range("N33").Value = -90602

# Row 47
$ws.Range("H47").Value = 2136.182
$ws.Range("I47").Value = 96
$ws.Range("J47").Value = 3836.3333
$ws.Range("K47").Value = 288
$ws.Range("L47").Value = 11508.9999
$ws.Range("M47").Value = 143
$ws.Range("N47").Value = -12370.9999

# Row 99
$ws.Range("H99").Value = 6399.1665
$ws.Range("I99").Value = 400
$ws.Range("J99").Value = 7599
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 22797
$ws.Range("M99").Value = 1046
$ws.Range("N99").Value = -27289

# Row 131
$ws.Range("H131").Value = 63334030
$ws.Range("I131").Value = 397
$ws.Range("J131").Value = 105556456
$ws.Range("K131").Value = 1191
$ws.Range("L131").Value = 316669368
$ws.Range("M131").Value = 3849
$ws.Range("N131").Value = -316679448

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 11667
$ws.Range("I12").Value = 11500.5
$ws.Range("J12").Value = 12000
$ws.Range("K12").Value = 11500.5
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = -11360.5
$ws.Range("N12").Value = -12280

# Row 132
$ws.Range("H132").Value = 58857.89
$ws.Range("I132").Value = 2282.2856
$ws.Range("J132").Value = 94860.55
$ws.Range("K132").Value = 6846.8568
$ws.Range("L132").Value = 284581.65
$ws.Range("M132").Value = -4316.8568
$ws.Range("N132").Value = -289641.65

$ws = $wb.Worksheets.Item("LTW")
# Row 47
$ws.Range("H47").Value = 49466.25
$ws.Range("J47").Value = 49466.25
$ws.Range("L47").Value = 49466.25
$ws.Range("N47").Value = -50446.25

# Row 52
$ws.Range("H52").Value = 49466.25
$ws.Range("J52").Value = 49466.25
$ws.Range("L52").Value = 49466.25
$ws.Range("N52").Value = -49932.25

# Row 132
$ws.Range("H132").Value = 525136.9
$ws.Range("I132").Value = 184200.45
$ws.Range("J132").Value = 837661.9399999999
$ws.Range("K132").Value = 552601.3500000001
$ws.Range("L132").Value = 2512985.82
$ws.Range("M132").Value = -550071.3500000001
$ws.Range("N132").Value = -2518045.82

$ws = $wb.Worksheets.Item("WVR")
# Row 28
$ws.Range("H28").Value = 48436.715
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

# Row 46
$ws.Range("H46").Value = 59990.5
$ws.Range("J46").Value = 59990.5
$ws.Range("L46").Value = 59990.5
$ws.Range("N46").Value = -60452.5

# Row 64
$ws.Range("H64").Value = 24995
$ws.Range("J64").Value = 24995
$ws.Range("L64").Value = 24995
$ws.Range("N64").Value = -25491

# Row 67
$ws.Range("H67").Value = 24995
$ws.Range("J67").Value = 24995
$ws.Range("L67").Value = 24995
$ws.Range("N67").Value = -26711

# Row 134
$ws.Range("H134").Value = 59990.5
$ws.Range("J134").Value = 59990.5
$ws.Range("L134").Value = 179971.5
$ws.Range("N134").Value = -185041.5
